$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Addresses" (sheet3): add helper column F that concatenates parts of
# the name/address columns, then flag duplicate results with conditional
# formatting ("Highlight Duplicate Values").
# ---------------------------------------------------------------------------
$wsAddr = $wb.Worksheets.Item("Addresses")

$wsAddr.Range("F1").Value = "Dupe Check"
$wsAddr.Range("F1").Font.Bold = $true

$wsAddr.Range("F2").Formula = "=LEFT(B2,3)&LEFT(C2,2)&MID(D2,2,4)&RIGHT(E2,3)"
$wsAddr.Range("F3").Formula = "=LEFT(B3,3)&LEFT(C3,2)&MID(D3,2,4)&RIGHT(E3,3)"
$wsAddr.Range("F3:F21").Formula = "=LEFT(B3,3)&LEFT(C3,2)&MID(D3,2,4)&RIGHT(E3,3)"

$fcAddr = $wsAddr.Range("F1:F1048576").FormatConditions.AddUniqueValues()
$fcAddr.DupeUnique = 1
$fcAddr.Font.Color = 393372
$fcAddr.Interior.Color = 13551615

$wsAddr.Range("F9").Select()

# ---------------------------------------------------------------------------
# Sheet "Names2" (sheet2): add helper column D that concatenates first/last
# name, then flag duplicates.
# ---------------------------------------------------------------------------
$wsNames2 = $wb.Worksheets.Item("Names2")

$wsNames2.Range("D1").Value = "Duplicate Check"
$wsNames2.Range("D1").Font.Bold = $true

$wsNames2.Range("D2").Formula = "=B2&C2"
$wsNames2.Range("D3").Formula = "=B3&C3"
$wsNames2.Range("D3:D38").Formula = "=B3&C3"

$fcNames2 = $wsNames2.Range("D1:D1048576").FormatConditions.AddUniqueValues()
$fcNames2.DupeUnique = 1
$fcNames2.Font.Color = 393372
$fcNames2.Interior.Color = 13551615

$wsNames2.Columns.Item(4).Select()

# ---------------------------------------------------------------------------
# Sheet "Names1" (sheet1): this is the cleaned-up roster, so remove the rows
# that are exact duplicates of another row, then flag any remaining
# duplicate student names with conditional formatting.
# ---------------------------------------------------------------------------
$wsNames1 = $wb.Worksheets.Item("Names1")

$rowsToDelete = @(33, 24, 20, 15, 14, 5, 3)
foreach ($r in $rowsToDelete) {
    $wsNames1.Rows.Item($r).Delete()
}

$wsNames1.Sort.SortFields.Clear()
$wsNames1.Sort.SortFields.Add($wsNames1.Range("B5:B31")) | Out-Null
$wsNames1.Sort.SetRange($wsNames1.Range("B2:B31"))
$wsNames1.Sort.Header = 2
$wsNames1.Sort.Apply()

$fcNames1 = $wsNames1.Range("B1:B1048576").FormatConditions.AddUniqueValues()
$fcNames1.DupeUnique = 1
$fcNames1.Font.Color = 393372
$fcNames1.Interior.Color = 13551615

$wsNames1.Range("C15").Select()

# Addresses tab ends up the active/selected sheet.
$wsAddr.Activate()
